$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: fill in payment date and collected amount (closes the 8110 balance)
$ws.Range("F7").Value = 44460
$ws.Range("G7").Value = 8110

# Row 8: new credit entry - OBRADOR
$ws.Range("A8").Value = 44460
$ws.Range("D5").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").WrapText = $true
$ws.Range("D8").Value = "OBRADOR"
$ws.Range("E8").Value = 604
$ws.Range("F8").Value = 44461
$ws.Range("G8").Value = 604

# Row 9: GUSTAVO
$ws.Range("A9").Value = 44462
$ws.Range("D9").Value = "GUSTAVO"
$ws.Range("E9").Value = 2028
$ws.Range("F9").Value = 44463
$ws.Range("G9").Value = 2028

# Row 10: GUSTAVO
$ws.Range("A10").Value = 44463
$ws.Range("D10").Value = "GUSTAVO"
$ws.Range("E10").Value = 1050
$ws.Range("F10").Value = 44464
$ws.Range("G10").Value = 1050

# Row 11: GUSTAVO
$ws.Range("A11").Value = 44463
$ws.Range("D11").Value = "GUSTAVO"
$ws.Range("E11").Value = 1786
$ws.Range("F11").Value = 44464
$ws.Range("G11").Value = 1786

# Row 12: EL PRIMO
$ws.Range("A12").Value = 44463
$ws.Range("D12").Value = "EL PRIMO"
$ws.Range("E12").Value = 420
$ws.Range("F12").Value = 44463
$ws.Range("G12").Value = 420

# Row 13: EL PRIMO
$ws.Range("A13").Value = 44463
$ws.Range("D13").Value = "EL PRIMO"
$ws.Range("E13").Value = 1284
$ws.Range("F13").Value = 44463
$ws.Range("G13").Value = 1284

# Update the active selection to reflect where the user left off editing
$ws.Range("E14").Select()
